$d = $word.ActiveDocument

# 1. Locate the discussion paragraph (it currently ends "...VPN routing
#    safely.") without hard-coding a paragraph index: search for a short,
#    unique phrase near its end and expand that hit to the whole
#    enclosing paragraph.
$hit = $d.Content.Duplicate
$hit.Find.Execute("VPN routing safely", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$hit.Expand(4) | Out-Null   # wdParagraph

# 2. Replace the lone trailing "." of that paragraph with the full
#    continuation sentence (itself still ending in "."). The trailing "."
#    is isolated as a 1-character Range immediately before the paragraph
#    mark, so the replacement lands exactly on that final run without
#    disturbing any other "." elsewhere in the document.
$periodRange = $d.Range($hit.End - 2, $hit.End - 1)
$periodRange.Text = ", the Hybrid solution of on premises and cloud storage is a good choice as the food chain organisation being large enough to maintain and handle costs so, it would viable for the organisation to choose a trusted vendor for SaaS and do the on premises implementation although there are economic factors this procedure helps to enhance productivity by allowing workers to work from anywhere and keep the intranet less likely to be attacked by implementing proper security management systems and network security such as TLS."

# 3. Give the trailing empty paragraph the same paragraph formatting as the
#    paragraph above it (1.5 line spacing, 0.5" left indent, justified).
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$lastPara.Format.LineSpacingRule = 1   # wdLineSpace1pt5
$lastPara.Format.LeftIndent = 36       # points (720 twips)
$lastPara.Format.Alignment = 3         # wdAlignParagraphJustify
